$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12; this shifts existing rows 12-55 down to 13-56
$ws.Rows("12").Insert()

# Populate the newly inserted row 12 with the new weekly price record
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C12").Value = "Arica y Parinacota"
$ws.Range("D12").Value = 45030
$ws.Range("E12").Value = 15
$ws.Range("F12").Value = "Fruta"
$ws.Range("G12").Value = 100108
$ws.Range("H12").Value = "Tropicales y subtropicales"
$ws.Range("I12").Value = 100108001
$ws.Range("J12").Value = "Guayaba"
$ws.Range("K12").Value = "Sin especificar"
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 100
$ws.Range("N12").Value = 1900
$ws.Range("O12").Value = 2000
$ws.Range("P12").Value = 1950
$ws.Range("Q12").Value = "`$/kilo (en caja de 10 kilos )"
$ws.Range("R12").Value = "Región de Arica y Parinacota"
$ws.Range("S12").Value = 1950
$ws.Range("T12").Value = 1
